$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fix "select" -> "dropdown" wording for Ticketart/Haeufigkeit question
$ws.Range("D4").Value = "Ticketart:dropdown(Bus,Zug,U-Bahn):pflicht;Häufigkeit:dropdown(Täglich,Wöchentlich,Selten):pflicht"

# Row 11: replace "Ja:boolean" attribute definition with the new checkbox question
$ws.Range("D11").Value = "Antworten komplett:checkbox(1,2,3):pflicht"

# Row 10: fix missing closing parenthesis typo in "Nicht-Dringend(5 Wochen und später"
$ws.Range("D10").Value = "Dringend(1 Woche):checkbox;Normal(2-4 Wochen):checkbox;Nicht-Dringend(5 Wochen und später):checkbox"

# Row 5: Upload answer changes from "Ja" to "Nein"
$ws.Range("F5").Value = "Nein"

# Row 4: Upload answer changes from "Nein" to literal text "TRUE" (kept as text, not boolean)
$ws.Range("F4").Value = "'TRUE"

# Update the last selected cell to reflect the author's final click position
$ws.Range("F5").Select()
